$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update cryptocurrency price and 1h-volume-change figures with latest scraped values.
# All target cells are plain text (inline strings) in the source sheet, so we force
# the text number format before assigning to prevent Excel from auto-coercing
# numeric-looking values (prices) into actual numbers.

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '56.573.02'
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = '  +4.32%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '3.008.20'
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = '  +5.08%  '
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.00'
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = '  +0.14%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '508.57'
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = '  +8.60%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '137.09'
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = '  +9.20%  '
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = '  -0.01%  '
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = '  +8.29%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '7.58'
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = '  +14.91%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.108'
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = '  +13.36%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.354'
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = '  +7.75%  '
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = '  +4.82%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '3.524.57'
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = '  +5.28%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '25.62'
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = '  +10.98%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.0000154'
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = '  +15.86%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '56.621.47'
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = '  +4.53%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '3.010.54'
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = '  +5.38%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '5.83'
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = '  +9.59%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '12.50'
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = '  +9.63%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '7.85'
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = '  +12.23%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '326.75'
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = '  +11.07%  '
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = '  -0.04%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '0.478'
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = '  +8.52%  '
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = '  +6.93%  '
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = '  +10.20%  '
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = '  +0.33%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '0.0₃0921'
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = '  +14.95%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '6.57'
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = '  +7.47%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '7.02'
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = '  +14.32%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '1.24'
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = '  +11.07%  '
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = '  +10.03%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '20.67'
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = '  +10.04%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '156.49'
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = '  +16.32%  '
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = '  +8.38%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '5.64'
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = '  +4.68%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '1.27'
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.0675'
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = '  +9.68%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '23.80'
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = '  +3.42%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '3.044.15'
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = '  +5.61%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '36.63'
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = '  +5.71%  '
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = '  +0.16%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.650'
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = '  +8.54%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '2.268.06'
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = '  +11.64%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '1.00'
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = '  +6.00%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '1.40'
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = '  +7.31%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '3.61'
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = '  +7.12%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '1.99'
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = '  +24.55%  '
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = '  +11.52%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '5.79'
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = '  +8.38%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '19.22'
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = '  +8.08%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.0875'
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = '  +11.28%  '
